$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear contents of A2:C4 and F2:J4, keeping styles on F:J columns
$ws.Range("A2:C4").ClearContents()
$ws.Range("F2:J4").ClearContents()

# Reset row heights to default for rows 2-4
$ws.Rows("2:4").AutoFit()
